# Children Tracking: "Switch map type" -> "Switch pick type"
# Adds two new localization rows (pick-up / drop-down switch labels) to
# the bottom of the Sheet1 translation table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting of the last "normal" data row (101) onto the two
# new rows so they pick up the same style (s="3" on all three columns).
$ws.Range("A101:C101").Copy()
$ws.Range("A104:C105").PasteSpecial(-4122)

# Fill in the new translation keys / values. The order in which new
# (not-yet-shared) strings are first entered controls their position in
# sharedStrings.xml, so we enter them in the same order the original
# author did: A104, B104, C104, B105, C105, A105.
$ws.Range("A104").Value = "lang_pick_type_UP"
$ws.Range("B104").Value = "Tuyến Đón"
$ws.Range("C104").Value = "Pick Up"

$ws.Range("B105").Value = "Tuyến Trả"
$ws.Range("C105").Value = "Drop Down"
$ws.Range("A105").Value = "lang_pick_type_DOWN"

# Scroll the view down to the new rows and select the last cell, matching
# the author's final cursor position.
$excel.ActiveWindow.ScrollRow = 91
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("A105").Select()
